$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.120.59"
$ws.Range("D3").Value = "1.654.74"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.93"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5248"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2607"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06353"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.37"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07804"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.504"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "1.641.83"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5477"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").Value = "0.0₅8206"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.40"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "26.116.92"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.577"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.63"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.06"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.030"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.92"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1244"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.257"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.15"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.429"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05914"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.254"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.588"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9532"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.790"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5699"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01619"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.788"
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8491"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.14"
$ws.Range("E42").Value = "  +2.94%  "
$ws.Range("D43").Value = "1.031.12"
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("D44").Value = "1.799.49"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.20"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9979"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.475"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05165"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.838"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09700"
$ws.Range("E51").Value = "  +0.18%  "
